{"js": "// Replace the two-digit multiplication problems in the document's table\n// with the new values from the commit diff. Each \"AxB=\" expression is\n// unique within the document, so a targeted search-and-replace for each\n// exact old string is safe and unambiguous.\nconst replacements = [\n  [\"23\u00d793=\", \"56\u00d749=\"],\n  [\"54\u00d764=\", \"53\u00d737=\"],\n  [\"59\u00d763=\", \"13\u00d742=\"],\n  [\"57\u00d768=\", \"71\u00d770=\"],\n  [\"92\u00d715=\", \"96\u00d719=\"],\n  [\"69\u00d753=\", \"60\u00d777=\"],\n  [\"46\u00d791=\", \"19\u00d782=\"],\n  [\"32\u00d735=\", \"46\u00d794=\"],\n  [\"45\u00d733=\", \"42\u00d788=\"],\n  [\"78\u00d767=\", \"49\u00d782=\"],\n  [\"69\u00d754=\", \"84\u00d717=\"],\n  [\"28\u00d722=\", \"59\u00d755=\"],\n  [\"77\u00d740=\", \"38\u00d746=\"],\n  [\"28\u00d741=\", \"86\u00d772=\"],\n  [\"31\u00d778=\", \"16\u00d748=\"],\n  [\"57\u00d785=\", \"11\u00d789=\"],\n  [\"19\u00d768=\", \"17\u00d768=\"],\n  [\"66\u00d759=\", \"41\u00d745=\"],\n  [\"79\u00d759=\", \"75\u00d764=\"],\n  [\"65\u00d739=\", \"48\u00d754=\"],\n  [\"34\u00d767=\", \"93\u00d741=\"],\n  [\"58\u00d717=\", \"61\u00d722=\"],\n  [\"92\u00d752=\", \"70\u00d765=\"],\n  [\"33\u00d783=\", \"13\u00d788=\"],\n  [\"81\u00d770=\", \"84\u00d789=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the two-digit multiplication problems in the document's table\n# with the new values from the commit diff. Each \"AxB=\" expression is\n# unique within the document, so a targeted Find/Replace for each exact\n# old string is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"23\u00d793=\", \"56\u00d749=\"),\n    @(\"54\u00d764=\", \"53\u00d737=\"),\n    @(\"59\u00d763=\", \"13\u00d742=\"),\n    @(\"57\u00d768=\", \"71\u00d770=\"),\n    @(\"92\u00d715=\", \"96\u00d719=\"),\n    @(\"69\u00d753=\", \"60\u00d777=\"),\n    @(\"46\u00d791=\", \"19\u00d782=\"),\n    @(\"32\u00d735=\", \"46\u00d794=\"),\n    @(\"45\u00d733=\", \"42\u00d788=\"),\n    @(\"78\u00d767=\", \"49\u00d782=\"),\n    @(\"69\u00d754=\", \"84\u00d717=\"),\n    @(\"28\u00d722=\", \"59\u00d755=\"),\n    @(\"77\u00d740=\", \"38\u00d746=\"),\n    @(\"28\u00d741=\", \"86\u00d772=\"),\n    @(\"31\u00d778=\", \"16\u00d748=\"),\n    @(\"57\u00d785=\", \"11\u00d789=\"),\n    @(\"19\u00d768=\", \"17\u00d768=\"),\n    @(\"66\u00d759=\", \"41\u00d745=\"),\n    @(\"79\u00d759=\", \"75\u00d764=\"),\n    @(\"65\u00d739=\", \"48\u00d754=\"),\n    @(\"34\u00d767=\", \"93\u00d741=\"),\n    @(\"58\u00d717=\", \"61\u00d722=\"),\n    @(\"92\u00d752=\", \"70\u00d765=\"),\n    @(\"33\u00d783=\", \"13\u00d788=\"),\n    @(\"81\u00d770=\", \"84\u00d789=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $true, $newText, 2) | Out-Null\n}\n"}
